# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they end up sharing the same bold/border/centered
# style used by the other headers, instead of Excel minting a brand new
# (duplicate) style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$IValues = @(1,7,9,6,7,5,5,6,2,2,5,5,5,4,3,1,1,1,1,5,6,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$JValues = @(3,8,9,9,7,6,7,6,6,4,7,8,7,5,6,5,6,5,6,7,7,5,6,5,2,7,6,6,7,6,6,5,5,4,3)

for ($i = 0; $i -lt $IValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $IValues[$i]
    $ws.Cells.Item($row, 10).Value = $JValues[$i]
}
